$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Export" sheet lists one Conta/Nome/Saldo row per account, sorted by
# balance descending. In the original layout, rows 4-6 were "ZULEIKA",
# "PEDRO" and "ERICA" (three high balances at the top of the list),
# followed by the big block "RODRIGO" .. "LUZIMAR" (rows 7-67, already
# sorted into place). The edit:
#   - removes "ZULEIKA" and "ERICA" entirely,
#   - keeps "PEDRO" but with a corrected/updated balance, and moves it
#     down so it lands in sorted order right after "LUZIMAR" (i.e. right
#     before "HEITOR", which used to be row 68).
# Net effect: the RODRIGO..LUZIMAR block shifts up to rows 4-64, "PEDRO"
# (balance 61.06) becomes row 65, and everything from "HEITOR" onward is
# unchanged save for shifting up by 2 rows overall.

# Stage PEDRO's account/name cells (A5:B5) out of the way before the rows
# around it are deleted, preserving their original text typing (so the
# leading zero in the account number survives the move intact) instead of
# re-typing them by hand.
$ws.Range("A5:B5").Copy($ws.Range("A400"))

# Drop the three original top rows: ZULEIKA (4), PEDRO (5), ERICA (6).
# This pulls the RODRIGO..LUZIMAR block up to rows 4-64.
$ws.Range("A4:C6").EntireRow.Delete()

# Make room for PEDRO again right after the block (now row 64, "LUZIMAR"),
# i.e. before what is now row 65 ("HEITOR"). The earlier delete shifted
# the staged cells from row 400 up to row 397, and this insert shifts them
# back down to row 398.
$ws.Range("A65:C65").EntireRow.Insert()
$ws.Range("A398:B398").Copy($ws.Range("A65"))
$ws.Range("C65").Value2 = 61.06

# Clean up the staging cells.
$ws.Range("A398:B398").ClearContents()
